# Asset Validation & Verification Portal - replace admin/security/verification
# test cases with the initial QR Code test cases (TC16, TC17, TC18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Multi-line "TEST STEPS" bodies (kept verbatim, including the double space
# after "1." in the first one and the trailing space at each line break).
# ---------------------------------------------------------------------------
$stepsTC16 = @'
1.  Login as Admin. 
2. Click "Templates" in the sidebar. 
3. Click "Reuse" your desired template. 
4. Upload excel file containing the certificate reciever's appropriate data (Name, email) 
5. Enable send Certificates.
6. Click Generate and Send.
7. Open your inbox and verify. (Also check your spam inbox)
'@

$stepsTC17 = @'
1. Login as Admin. 
2. Click "Templates" in the sidebar. 
3. Click "Reuse" your desired template. 
4. Upload excel file containing the certificate reciever's appropriate data (Name, email) 
5. Enable send Certificates. 
6. Click Generate and Send. 
7. Open your inbox and verify. (Also check your spam inbox)
'@

# ---------------------------------------------------------------------------
# Row 2 -> TC16 / QR Code
# (column order A, C, D, F, E mirrors how the author originally typed the
# row, which is what determines shared-string table ordering on save)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "TC16"
$ws.Range("C2").Value = "QR Code"
$ws.Range("D2").Value = "Verify QR code is clearly visible and not blurred on the generated file."
$ws.Range("F2").Value = "The QR Code is visible and clear."
$ws.Range("E2").Value = $stepsTC16

# ---------------------------------------------------------------------------
# Row 4 -> TC17 / QR Code
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "TC17"
$ws.Range("C4").Value = "QR Code"
$ws.Range("D4").Value = "Verify QR code scans correctly using a standard mobile camera app."
$ws.Range("F4").Value = "The QR Code is scannable by a mobile camera app."
$ws.Range("E4").Value = $stepsTC17

# ---------------------------------------------------------------------------
# Row 6 -> TC18 / QR Code
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "TC18"
$ws.Range("C6").Value = "QR Code"
$ws.Range("D6").Value = "Verify QR code redirects to the correct unique verification URL."
$ws.Range("F6").Value = "The QR Code redirects to the verification Portal of CertifyHub."
$ws.Range("E6").Value = $stepsTC17

# ---------------------------------------------------------------------------
# Drop the remaining old test cases (former rows 7-14: TC21/TC24/TC25/TC32).
# ---------------------------------------------------------------------------
$ws.Rows("7:14").Delete()

# ---------------------------------------------------------------------------
# Row heights for the three remaining data rows (grew to fit the new,
# longer wrapped "TEST STEPS" text).
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 106.8
$ws.Rows.Item(4).RowHeight = 93.6
$ws.Rows.Item(6).RowHeight = 93.6

# ---------------------------------------------------------------------------
# Column widths (re-tuned by the author alongside the content edit).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.833333333333334
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 16
$ws.Columns.Item(4).ColumnWidth = 38.166666666666664
$ws.Columns.Item(5).ColumnWidth = 75.83333333333333
$ws.Columns.Item(6).ColumnWidth = 39.5
$ws.Columns.Item(7).ColumnWidth = 33

# ---------------------------------------------------------------------------
# View: zoom out to 70% and move the active selection to F18.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("F18").Select()
